# Update the "Förändrad" (Changed) date column C for all data rows (2-388)
# from serial date 46060 (2026-02-07) to 46061 (2026-02-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 388
$range = $ws.Range("C2:C$lastRow")
$range.Value = 46061
